# "13 11 2021 - part 2"
# Adds new kas-kecil transactions (13 Nov 2021) to the "November" sheet:
#   - marks row 23 ("uang masuk ruko rich palace...") as "lunas" in col K
#   - clears the now-unused SUM(D23:D33) helper formula in J24
#   - fills in three new rows (38-40) of income/expense entries
#   - updates the current selection to E26

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November")
$ws.Activate()

# Row 23: mark this entry as settled ("lunas") in the new K column.
$ws.Range("K23").Value2 = "lunas"

# Row 24: the SUM(D23:D33) subtotal formula is no longer needed - clear it
# but keep the existing cell formatting/style.
$ws.Range("J24").ClearContents()

# New date cells (B38:B40) need the same date number format as the
# surrounding entries, so copy B37's formatting down first, then fill in
# the actual values/formulas for each row.
$ws.Range("B37").Copy($ws.Range("B38"))
$ws.Range("B37").Copy($ws.Range("B39"))
$ws.Range("B37").Copy($ws.Range("B40"))

# Row 38: uang masuk kas kecil asrama dari ce nanda (1.250.000)
$ws.Range("B38").Value2 = 44513
$ws.Range("C38").Value2 = 1250000
$ws.Range("D38").Value2 = 0
$ws.Range("F38").Value2 = "uang masuk kas kecil asrama dari ce nanda"
$ws.Range("G38").Value2 = "yofandi"

# Row 39: uang beli barang keperluan asrama di indogrosir (1.576.500)
$ws.Range("B39").Value2 = 44513
$ws.Range("C39").Value2 = 0
$ws.Range("D39").Value2 = 1576500
$ws.Range("F39").Value2 = "uang beli barang keperluan asrama di indogrosir"
$ws.Range("G39").Value2 = "yofandi sama valen sama jhonan"

# Row 40: uang beli bensin motor pertalite (15.000)
$ws.Range("B40").Value2 = 44513
$ws.Range("C40").Value2 = 0
$ws.Range("D40").Value2 = 15000
$ws.Range("F40").Value2 = "uang beli bensin motor pertalite"
$ws.Range("G40").Value2 = "saferius sama peter"

# Widen column G a bit to fit the new longer notes.
$ws.Columns.Item(7).ColumnWidth = 30.7

# Leave the same cell selected as in the saved file.
$ws.Range("E26").Select() | Out-Null
